$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Preserve the existing style while forcing the value to be stored as text,
    # so purely-numeric-looking strings (e.g. "116.62") are not coerced into numbers.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "51.663.58"
$ws.Range("E2").Value = "  +4.44%  "

$ws.Range("D3").Value = "2.766.80"
$ws.Range("E3").Value = "  +5.28%  "

$ws.Range("E4").Value = "  -0.02%  "

Set-TextCell $ws.Range("D5") "116.62"
$ws.Range("E5").Value = "  +3.97%  "

Set-TextCell $ws.Range("D6") "333.25"
$ws.Range("E6").Value = "  +2.93%  "

Set-TextCell $ws.Range("D7") "0.539"
$ws.Range("E7").Value = "  +2.68%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +6.24%  "

Set-TextCell $ws.Range("D10") "42.07"
$ws.Range("E10").Value = "  +6.31%  "

Set-TextCell $ws.Range("D11") "0.0861"
$ws.Range("E11").Value = "  +6.29%  "

Set-TextCell $ws.Range("D12") "20.29"
$ws.Range("E12").Value = "  +2.85%  "

$ws.Range("E13").Value = "  +2.22%  "

Set-TextCell $ws.Range("D14") "7.66"
$ws.Range("E14").Value = "  +5.26%  "

$ws.Range("D15").Value = "3.206.11"
$ws.Range("E15").Value = "  +7.19%  "

$ws.Range("D16").Value = "2.761.19"
$ws.Range("E16").Value = "  +5.22%  "

$ws.Range("E17").Value = "  +4.38%  "

$ws.Range("D18").Value = "51.642.74"
$ws.Range("E18").Value = "  +4.49%  "

$ws.Range("E19").Value = "  +13.42%  "

$ws.Range("E20").Value = "  +5.08%  "

$ws.Range("E21").Value = "  +2.89%  "

$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  +3.22%  "

Set-TextCell $ws.Range("D23") "278.30"
$ws.Range("E23").Value = "  +3.40%  "

Set-TextCell $ws.Range("D24") "69.89"
$ws.Range("E24").Value = "  +1.37%  "

$ws.Range("E25").Value = "  +6.22%  "

Set-TextCell $ws.Range("D26") "26.83"
$ws.Range("E26").Value = "  +2.35%  "

$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("E30").Value = "  +3.73%  "

Set-TextCell $ws.Range("D31") "35.05"
$ws.Range("E31").Value = "  +1.16%  "

Set-TextCell $ws.Range("D32") "50.17"
$ws.Range("E32").Value = "  +1.58%  "

$ws.Range("E33").Value = "  +1.54%  "

Set-TextCell $ws.Range("D34") "0.0820"
$ws.Range("E34").Value = "  +0.88%  "

Set-TextCell $ws.Range("D35") "19.20"
$ws.Range("E35").Value = "  +2.19%  "

$ws.Range("E36").Value = "  -0.16%  "

Set-TextCell $ws.Range("D37") "5.03"
$ws.Range("E37").Value = "  +3.22%  "

Set-TextCell $ws.Range("D38") "2.09"
$ws.Range("E38").Value = "  +2.87%  "

$ws.Range("E39").Value = "  +4.64%  "

Set-TextCell $ws.Range("D40") "0.0355"
$ws.Range("E40").Value = "  +9.55%  "

Set-TextCell $ws.Range("D41") "128.28"
$ws.Range("E41").Value = "  -0.27%  "

Set-TextCell $ws.Range("D42") "23.29"
$ws.Range("E42").Value = "  +5.03%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D43") "0.114"
$ws.Range("E43").Value = "  +3.61%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws.Range("D44") "2.31"
$ws.Range("E44").Value = "  +7.69%  "

Set-TextCell $ws.Range("D45") "2.45"
$ws.Range("E45").Value = "  +16.94%  "

$ws.Range("D46").Value = "2.088.93"
$ws.Range("E46").Value = "  +1.80%  "

$ws.Range("E48").Value = "  +3.82%  "

Set-TextCell $ws.Range("D49") "5.54"
$ws.Range("E49").Value = "  +6.85%  "

Set-TextCell $ws.Range("D50") "60.44"
$ws.Range("E50").Value = "  +3.04%  "

$ws.Range("E51").Value = "  -0.45%  "
